# Issue 43 (csarven): fix "qb:dataset" -> "qb:dataSet" label and widen its textbox.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item("TextBox 30")

$shp.TextFrame.TextRange.Text = "qb:dataSet"
$shp.Width = 68.1845
